$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menu")

# Style tweak: B6 and B7 get the light grey highlight style (same as B3/B5/B11-B17...)
$ws.Range("B3").Copy()
$ws.Range("B6:B7").PasteSpecial(-4122)

# Insert a new row above the old row 26 ("Asignación" block), shifting everything below down by one.
$ws.Rows.Item(25).Insert()
$ws.Rows.Item(25).RowHeight = 13.8

# New row content: a "Relacion Proveedor - Producto" menu entry in column B.
$ws.Range("B25").Value = "Relacion Proveedor " + [char]0x2013 + " Producto"

# Match the header-row look (bold + grey fill) used elsewhere in row 2.
$ws.Range("A2").Copy()
$ws.Range("B25").PasteSpecial(-4122)

# Restore the saved selection/active-cell state.
$ws.Activate()
$ws.Range("A13").Select()
